# natmiOut/OldD0/LR-pairs_lrc2p/Pgf-Nrp2.xlsx -- "Natmi following Dr Hou advice"
#
# The clustering used to produce this ligand(Pgf)-receptor(Nrp2) edge table was
# rerun with an extra sending/target cluster (ECs), so the cluster set is now
# {ECs, FAPs, sCs} instead of just {FAPs, sCs}. That makes the sending x target
# combination grid grow from 2x3=6 rows to 3x3=9 rows (sheet rows 2-10 instead of
# 2-7), and every specificity/expression statistic is recomputed for the new run.
# Rewrite the whole data block (A2:T10) with the refreshed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Pgf"
$ws.Range("C2").Value = "Nrp2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 22.51188433333333
$ws.Range("H2").Value = 67.535653
$ws.Range("I2").Value = 0.7173237801266834
$ws.Range("J2").Value = 0.7173237801266834
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 26.83081766666667
$ws.Range("N2").Value = 80.49245300000001
$ws.Range("O2").Value = 0.5916656861001716
$ws.Range("P2").Value = 0.5916656861001716
$ws.Range("Q2").Value = 604.0122638807566
$ws.Range("R2").Value = 5436.11037492681
$ws.Range("S2").Value = 0.4244158665246228
$ws.Range("T2").Value = 0.4244158665246228

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Pgf"
$ws.Range("C3").Value = "Nrp2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 22.51188433333333
$ws.Range("H3").Value = 67.535653
$ws.Range("I3").Value = 0.7173237801266834
$ws.Range("J3").Value = 0.7173237801266834
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.21969166666667
$ws.Range("N3").Value = 30.659075
$ws.Range("O3").Value = 0.2253617819930474
$ws.Range("P3").Value = 0.2253617819930474
$ws.Range("Q3").Value = 230.0645167223306
$ws.Range("R3").Value = 2070.580650500975
$ws.Range("S3").Value = 0.1616573653553383
$ws.Range("T3").Value = 0.1616573653553383

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Pgf"
$ws.Range("C4").Value = "Nrp2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 22.51188433333333
$ws.Range("H4").Value = 67.535653
$ws.Range("I4").Value = 0.7173237801266834
$ws.Range("J4").Value = 0.7173237801266834
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 8.297426666666667
$ws.Range("N4").Value = 24.89228
$ws.Range("O4").Value = 0.1829725319067811
$ws.Range("P4").Value = 0.1829725319067811
$ws.Range("Q4").Value = 186.7907093843155
$ws.Range("R4").Value = 1681.11638445884
$ws.Range("S4").Value = 0.1312505482467224
$ws.Range("T4").Value = 0.1312505482467224

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Pgf"
$ws.Range("C5").Value = "Nrp2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 7.074492666666667
$ws.Range("H5").Value = 21.223478
$ws.Range("I5").Value = 0.2254232363222357
$ws.Range("J5").Value = 0.2254232363222357
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 26.83081766666667
$ws.Range("N5").Value = 80.49245300000001
$ws.Range("O5").Value = 0.5916656861001716
$ws.Range("P5").Value = 0.5916656861001716
$ws.Range("Q5").Value = 189.8144228235038
$ws.Range("R5").Value = 1708.329805411534
$ws.Range("S5").Value = 0.1333751937815167
$ws.Range("T5").Value = 0.1333751937815167

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Pgf"
$ws.Range("C6").Value = "Nrp2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 7.074492666666667
$ws.Range("H6").Value = 21.223478
$ws.Range("I6").Value = 0.2254232363222357
$ws.Range("J6").Value = 0.2254232363222357
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 10.21969166666667
$ws.Range("N6").Value = 30.659075
$ws.Range("O6").Value = 0.2253617819930474
$ws.Range("P6").Value = 0.2253617819930474
$ws.Range("Q6").Value = 72.29913375142779
$ws.Range("R6").Value = 650.69220376285
$ws.Range("S6").Value = 0.05080178224021888
$ws.Range("T6").Value = 0.05080178224021888

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Pgf"
$ws.Range("C7").Value = "Nrp2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 7.074492666666667
$ws.Range("H7").Value = 21.223478
$ws.Range("I7").Value = 0.2254232363222357
$ws.Range("J7").Value = 0.2254232363222357
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 8.297426666666667
$ws.Range("N7").Value = 24.89228
$ws.Range("O7").Value = 0.1829725319067811
$ws.Range("P7").Value = 0.1829725319067811
$ws.Range("Q7").Value = 58.70008410553778
$ws.Range("R7").Value = 528.30075694984
$ws.Range("S7").Value = 0.04124626030050012
$ws.Range("T7").Value = 0.04124626030050011

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Pgf"
$ws.Range("C8").Value = "Nrp2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.796779333333333
$ws.Range("H8").Value = 5.390338
$ws.Range("I8").Value = 0.05725298355108089
$ws.Range("J8").Value = 0.05725298355108089
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 26.83081766666667
$ws.Range("N8").Value = 80.49245300000001
$ws.Range("O8").Value = 0.5916656861001716
$ws.Range("P8").Value = 0.5916656861001716
$ws.Range("Q8").Value = 48.20905867990157
$ws.Range("R8").Value = 433.8815281191141
$ws.Range("S8").Value = 0.03387462579403211
$ws.Range("T8").Value = 0.03387462579403211

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Pgf"
$ws.Range("C9").Value = "Nrp2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.796779333333333
$ws.Range("H9").Value = 5.390338
$ws.Range("I9").Value = 0.05725298355108089
$ws.Range("J9").Value = 0.05725298355108089
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 10.21969166666667
$ws.Range("N9").Value = 30.659075
$ws.Range("O9").Value = 0.2253617819930474
$ws.Range("P9").Value = 0.2253617819930474
$ws.Range("Q9").Value = 18.36253077970556
$ws.Range("R9").Value = 165.26277701735
$ws.Range("S9").Value = 0.01290263439749022
$ws.Range("T9").Value = 0.01290263439749022

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Pgf"
$ws.Range("C10").Value = "Nrp2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.796779333333333
$ws.Range("H10").Value = 5.390338
$ws.Range("I10").Value = 0.05725298355108089
$ws.Range("J10").Value = 0.05725298355108089
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 8.297426666666667
$ws.Range("N10").Value = 24.89228
$ws.Range("O10").Value = 0.1829725319067811
$ws.Range("P10").Value = 0.1829725319067811
$ws.Range("Q10").Value = 14.90864475451555
$ws.Range("R10").Value = 134.17780279064
$ws.Range("S10").Value = 0.01047572335955856
$ws.Range("T10").Value = 0.01047572335955856
